$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$newLine = "ChangeConfigxml(Configuration,WebServer,<WebServer>endl  <Enabled VALUE=`"1`"/>endl  <Port VALUE=`"8082`"/>endl  <WebFolder VALUE=`"\\auto\\ComplianceTest_JS\`"/>endl  <Public VALUE=`"1`"/>endl</WebServer>endl);"

$text = "wait(3);`nPullConfigxml;`nChangeConfigxml(Configuration/Applications/Application/General,StartPage,<StartPage value=`"http://127.0.0.1:8082/app/`" name=`"Menu`"/>);`n" + $newLine + "`nChangeConfigxml(Configuration/Screen,FullScreen,<FullScreen value=`"0`"/>);`nPushConfigxml;"

$ws.Range("G2").Value = $text

$ws.Rows.Item(2).RowHeight = 243
